$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''39.754.25'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '''2.187.65'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.12%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''293.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").Value = '''86.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("E7").Value = '  -1.78%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -1.50%  '
$ws.Range("D10").Value = '''29.82'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.85%  '
$ws.Range("D11").Value = '''0.0774'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.27%  '
$ws.Range("D12").Value = '''49.65'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.69%  '
$ws.Range("E13").Value = '  +2.34%  '
$ws.Range("D14").Value = '''6.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("D15").Value = '''2.528.68'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.13%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '''13.59'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.84%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '''2.104.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.40%  '
$ws.Range("D18").Value = '''0.720'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.53%  '
$ws.Range("D19").Value = '''39.679.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.80%  '
$ws.Range("D20").Value = '''0.0₃0877'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.67%  '
$ws.Range("D21").Value = '''11.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.42%  '
$ws.Range("D22").Value = '''5.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.00%  '
$ws.Range("D23").Value = '''64.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.94%  '
$ws.Range("D24").Value = '''235.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("E26").Value = '  -1.31%  '
$ws.Range("E27").Value = '  -3.12%  '
$ws.Range("D28").Value = '''22.35'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.18%  '
$ws.Range("E29").Value = '  -3.59%  '
$ws.Range("D30").Value = '''9.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.57%  '
$ws.Range("D31").Value = '''156.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.36%  '
$ws.Range("D32").Value = '''31.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.77%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").Value = '''4.85'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.86%  '
$ws.Range("D35").Value = '''0.0703'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.56%  '
$ws.Range("E36").Value = '  -2.60%  '
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("D38").Value = '''2.78'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.42%  '
$ws.Range("D39").Value = '''0.0966'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.02%  '
$ws.Range("D40").Value = '''15.14'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.04%  '
$ws.Range("E41").Value = '  -3.08%  '
$ws.Range("D42").Value = '''2.106.10'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.72%  '
$ws.Range("E43").Value = '  -3.06%  '
$ws.Range("E44").Value = '  -2.03%  '
$ws.Range("E45").Value = '  -2.09%  '
$ws.Range("D46").Value = '''17.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.68%  '
$ws.Range("D47").Value = '''9.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.41%  '
$ws.Range("D48").Value = '''2.63'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.74%  '
$ws.Range("D49").Value = '''2.401.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.76%  '
$ws.Range("E50").Value = '  +1.75%  '
$ws.Range("E51").Value = '  +0.29%  '
